# Regenerate merged AHB files
#
# The sheet holds a "before/after" diff table (AHB-Diff) whose header row
# used generic "_old" / "_new" suffixes. This edit re-labels the headers with
# the concrete format-version identifiers that were actually diffed
# (FV2304 -> FV2310), turns the data range into a proper Excel Table so the
# generated file is filterable/sortable, and freezes the header row so it
# stays visible while scrolling through the 58 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Rename the header-row labels: "<Field>_old" -> "<Field>_FV2304" and
#    "<Field>_new" -> "<Field>_FV2310". Column K ("diff") is left as-is.
# ---------------------------------------------------------------------
$fieldNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $fieldNames.Length; $i++) {
    # Columns A..J (1..10) -> "_FV2304"
    $ws.Cells.Item(1, $i + 1).Value = $fieldNames[$i] + "_FV2304"
    # Columns L..U (12..21) -> "_FV2310"
    $ws.Cells.Item(1, $i + 12).Value = $fieldNames[$i] + "_FV2310"
}

# ---------------------------------------------------------------------
# 2. Freeze the header row (split under row 1, frozen, active pane bottom
#    left) so headers remain visible while scrolling.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------
# 3. Turn the whole populated range into an Excel Table ("Table1") with
#    an AutoFilter, matching the header row exactly.
# ---------------------------------------------------------------------
$dataRange = $ws.Range("A1:U58")
$lo = $ws.ListObjects.Add(1, $dataRange, 0, 1)
$lo.Name = "Table1"
